# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect a
# completed handback: the Overview status text is refreshed, the
# "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns on the per-language sheets are populated with the
# handback xliff file names and timestamps (with working hyperlinks to
# the source markdown on the Target File column, matching column A),
# and a couple of columns are widened so the new, longer values are
# readable.

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8e0fc4f515407a7a183650094c4955c2f74fd203/e2e/"
$mdFile1 = "6c864454-4b85-46d2-be32-8d1575f62b92.md"
$mdFile2 = "9c93726a-a17b-4071-91f6-485915fec7c8.md"

# ---------------------------------------------------------------------
# 1) Overview sheet: the handoff is now a handback that is in sync with
#    en-US, and the (now longer) status text needs wider columns.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"

$overview.Columns("E:F").ColumnWidth = 29.144371396019366

# ---------------------------------------------------------------------
# 2) Per-language sheets: zh-cn (sheet 2) and de-de (sheet 3) each get
#    the Latest Target File, Latest Handback File and Latest Handback
#    DateTime columns (I, J, K) filled in for both data rows.
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# zh-cn -----------------------------------------------------------------
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), ($repoBase + $mdFile1), "", "", $mdFile1)
$zhcn.Range("J2").Value = "6c864454-4b85-46d2-be32-8d1575f62b92.be0966b3a36dfc5d82d17d618880c3c24567f6f2.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-17 21:02:31"

$zhcn.Hyperlinks.Add($zhcn.Range("I3"), ($repoBase + $mdFile2), "", "", $mdFile2)
$zhcn.Range("J3").Value = "9c93726a-a17b-4071-91f6-485915fec7c8.c71666dde1fbf2c7e6c1ab9f67e2953b9dc571e5.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-17 21:02:31"

$zhcn.Columns("C:C").ColumnWidth = 29.144371396019366
$zhcn.Columns("I:J").ColumnWidth = 39.16666666666666

# de-de -------------------------------------------------------------------
$dede.Hyperlinks.Add($dede.Range("I2"), ($repoBase + $mdFile1), "", "", $mdFile1)
$dede.Range("J2").Value = "6c864454-4b85-46d2-be32-8d1575f62b92.be0966b3a36dfc5d82d17d618880c3c24567f6f2.de-de.xlf"
$dede.Range("K2").Value = "2016-08-17 21:02:38"

$dede.Hyperlinks.Add($dede.Range("I3"), ($repoBase + $mdFile2), "", "", $mdFile2)
$dede.Range("J3").Value = "9c93726a-a17b-4071-91f6-485915fec7c8.c71666dde1fbf2c7e6c1ab9f67e2953b9dc571e5.de-de.xlf"
$dede.Range("K3").Value = "2016-08-17 21:02:38"

$dede.Columns("C:C").ColumnWidth = 29.144371396019366
$dede.Columns("I:J").ColumnWidth = 39.16666666666666

Write-Host "Handback report generated."
